# Re-upload of yahoo/NSY.xlsx: the player/position/team table (A2:C17) is
# reordered, and Dennis Schröder's team is updated from the Brooklyn Nets
# to the Golden State Warriors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("T.J. McConnell",           "PG",         "Indiana Pacers"),
    @("Jordan Poole",             "PG,SG",      "Washington Wizards"),
    @("Mike Conley",              "PG",         "Minnesota Timberwolves"),
    @("Kyrie Irving",             "PG,SG",      "Dallas Mavericks"),
    @("CJ McCollum",              "PG,SG",      "New Orleans Pelicans"),
    @("Dennis Schröder",          "PG",         "Golden State Warriors"),
    @("Lauri Markkanen",          "SF,PF",      "Utah Jazz"),
    @("Jimmy Butler",             "SF,PF",      "Miami Heat"),
    @("Zach LaVine",              "SG,SF",      "Chicago Bulls"),
    @("RJ Barrett",               "SF,PF",      "Toronto Raptors"),
    @("Joel Embiid",              "C",          "Philadelphia 76ers"),
    @("John Collins",             "PF,C",       "Utah Jazz"),
    @("Jalen Williams",           "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Shai Gilgeous-Alexander",  "PG",         "Oklahoma City Thunder"),
    @("Keyonte George",           "PG,SG",      "Utah Jazz"),
    @("Tobias Harris",            "SF,PF",      "Detroit Pistons")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
